$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row: "<Name>_old" -> "<Name>_FV2304", "<Name>_new" ->
#    "<Name>_FV2310". Column K ("diff") is left untouched.
# ---------------------------------------------------------------------------
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $header = $cell.Value()
    if ($header -like "*_old") {
        $cell.Value = ($header -replace "_old$", "_FV2304")
    } elseif ($header -like "*_new") {
        $cell.Value = ($header -replace "_new$", "_FV2310")
    }
}

# ---------------------------------------------------------------------------
# 2) Turn the data range into a native Excel Table ("Table1").
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:U93")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split below row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
